$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.76329235263777306
$ws.Range("C1").Value = 0.8841818121120798
$ws.Range("BO1").Value = 0.92846320351391554
$ws.Range("BP1").Value = 0.96594238092680706
$ws.Range("C2").Value = 0.95739534546396388
$ws.Range("M2").Value = 0.82972190159711801
$ws.Range("P2").Value = 0.52430751820335753
$ws.Range("D3").Value = 0.94493407009390951
$ws.Range("E3").Value = 0.84929483246677251
$ws.Range("B4").Value = 0.93294352633416699
$ws.Range("D5").Value = 0.85863121310145607
$ws.Range("AY5").Value = 0.86617455020973932
$ws.Range("D6").Value = 0.97985703602802698
$ws.Range("E6").Value = 0.89323393204577939
$ws.Range("H6").Value = 0.7768515609129365
$ws.Range("AX6").Value = 0.85084054636879869
$ws.Range("H7").Value = 0.64154393302573443
$ws.Range("I7").Value = 0.95868352517013478
$ws.Range("I8").Value = 0.80773496365022401
$ws.Range("J8").Value = 0.64155835827222707
$ws.Range("J9").Value = 0.82744893414250364
$ws.Range("K10").Value = 0.59413475611487521
$ws.Range("O10").Value = 0.94196035399956424
$ws.Range("I11").Value = 0.9765218950849881
$ws.Range("AV11").Value = 0.90567464974268597
$ws.Range("J12").Value = 0.53756482940775441
$ws.Range("K12").Value = 0.91217199882082123
$ws.Range("K13").Value = 0.8256361917940036
$ws.Range("L13").Value = 0.97895677646059931
$ws.Range("O13").Value = 0.74941848019956558
$ws.Range("L14").Value = 0.76189824052786992
$ws.Range("AE14").Value = 0.66569101259080643
$ws.Range("N15").Value = 0.59004007128098501
$ws.Range("Q15").Value = 0.6722929629587362
$ws.Range("Q16").Value = 0.78822577747232614
$ws.Range("R17").Value = 0.81014838557996915
$ws.Range("S18").Value = 0.66444500801071138
$ws.Range("T18").Value = 0.76052017887150569
$ws.Range("Q19").Value = 0.67995909487666206
$ws.Range("AB19").Value = 0.56376741449053669
$ws.Range("S20").Value = 0.76723183398878469
$ws.Range("U20").Value = 0.69316603606553462
$ws.Range("AD21").Value = 0.95360801465817691
$ws.Range("T22").Value = 0.75723572278312323
$ws.Range("U22").Value = 0.81170673906379076
$ws.Range("X22").Value = 0.92313827115283753
$ws.Range("V23").Value = 0.88130917683759469
$ws.Range("X23").Value = 0.96321852834708421
$ws.Range("Y23").Value = 0.91024000674990335
$ws.Range("Y24").Value = 0.94150737090626735
$ws.Range("AT24").Value = 0.8339947394556354
$ws.Range("H25").Value = 0.64459168828583624
$ws.Range("E26").Value = 0.6319497866327386
$ws.Range("Y26").Value = 0.98963834513767224
$ws.Range("AD26").Value = 0.95685793425911003
$ws.Range("AE26").Value = 0.79465621197834846
$ws.Range("AB27").Value = 0.67675048239929203
$ws.Range("AD28").Value = 0.97192544231147271
$ws.Range("AA29").Value = 0.56835634227137422
$ws.Range("AB29").Value = 0.85399051090348366
$ws.Range("AG31").Value = 0.99253572654051325
$ws.Range("AM31").Value = 0.80748478973177762
$ws.Range("AD32").Value = 0.97736451615776976
$ws.Range("AE32").Value = 0.87813315170379913
$ws.Range("Q33").Value = 0.77382436349056549
$ws.Range("AF33").Value = 0.99635278369879487
$ws.Range("AJ34").Value = 0.9214196053876359
$ws.Range("AG35").Value = 0.69449811034684661
$ws.Range("AH35").Value = 0.84244304822298177
$ws.Range("AK35").Value = 0.79220720887830298
$ws.Range("AI36").Value = 0.74701299737812121
$ws.Range("AL37").Value = 0.78929308735216552
$ws.Range("AM37").Value = 0.67789193930786062
$ws.Range("AJ38").Value = 0.73542610467774838
$ws.Range("AL39").Value = 0.55159239077818467
$ws.Range("BP39").Value = 0.95092576587190747
$ws.Range("AL40").Value = 0.60326531075446055
$ws.Range("AO40").Value = 0.99634779028825848
$ws.Range("AQ41").Value = 0.77543284414453728
$ws.Range("BJ41").Value = 0.82458088913237781
$ws.Range("AN42").Value = 0.87367329491451828
$ws.Range("AQ42").Value = 0.80378725391597294
$ws.Range("BP43").Value = 0.51648473272991302
$ws.Range("AP44").Value = 0.84219074138508865
$ws.Range("AQ44").Value = 0.83651421373909862
$ws.Range("AF45").Value = 0.64102489874888224
$ws.Range("AQ45").Value = 0.74257659640098383
$ws.Range("AR45").Value = 0.81384295396953443
$ws.Range("AT45").Value = 0.97493104806614383
$ws.Range("AC46").Value = 0.51588690541061299
$ws.Range("AR46").Value = 0.90895363320581102
$ws.Range("AU46").Value = 0.85693900175396465
$ws.Range("AV46").Value = 0.59168156296855234
$ws.Range("AJ47").Value = 0.92091818513214951
$ws.Range("AP47").Value = 0.81922134837438665
$ws.Range("AS47").Value = 0.8279457001194428
$ws.Range("AX48").Value = 0.81139128323307874
$ws.Range("BE48").Value = 0.85051296054255188
$ws.Range("L49").Value = 0.83402002443033374
$ws.Range("AW50").Value = 0.93125028431753676
$ws.Range("BP50").Value = 0.88135706118260104
$ws.Range("BA51").Value = 0.8733547613313567
$ws.Range("AJ52").Value = 0.58396967656601118
$ws.Range("AX52").Value = 0.69701790450589707
$ws.Range("AY52").Value = 0.85662247605479436
$ws.Range("BA52").Value = 0.58201369263758862
$ws.Range("AC53").Value = 0.73296313692616821
$ws.Range("BB53").Value = 0.9749300549719413
$ws.Range("AS54").Value = 0.85384687157262951
$ws.Range("AZ54").Value = 0.82975498986624707
$ws.Range("BD55").Value = 0.90991862578808402
$ws.Range("BB56").Value = 0.69918009734045961
$ws.Range("BD57").Value = 0.86331496292178023
$ws.Range("BF57").Value = 0.98034333026071607
$ws.Range("BD58").Value = 0.86882912645722499
$ws.Range("BG58").Value = 0.61768781133396833
$ws.Range("BH58").Value = 0.5786135803367991
$ws.Range("BH59").Value = 0.725160372851271
$ws.Range("BI60").Value = 0.68425254397831914
$ws.Range("AE61").Value = 0.9394639674340125
$ws.Range("BG61").Value = 0.93413405636628433
$ws.Range("AB62").Value = 0.78724132372435851
$ws.Range("BH62").Value = 0.88006882635173689
$ws.Range("BI62").Value = 0.79399185284975438
$ws.Range("BJ63").Value = 0.61673711152228872
$ws.Range("BK64").Value = 0.73264555299803336
$ws.Range("BC65").Value = 0.98341177896676957
$ws.Range("BI65").Value = 0.69286872151588086
$ws.Range("BK65").Value = 0.85601819544086244
$ws.Range("BL65").Value = 0.99082206927518524
$ws.Range("BN65").Value = 0.94811100423946704
$ws.Range("AC66").Value = 0.95728799342990512
$ws.Range("BL66").Value = 0.69714906126543319
$ws.Range("BO66").Value = 0.66325877734840777
$ws.Range("BP66").Value = 0.81587847531488455
$ws.Range("BP67").Value = 0.79698805270559037
